# Refresh the coinranking price/volume/coin snapshot to match the new scrape
# (commit: "Updated symbol list on Wed Jan 18 12:50:59 UTC 2023 with GitHub Actions").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns D (Price) and E (Volume(1h)) hold numeric-looking text (e.g. "301.59",
# "-0.20%"); a leading apostrophe forces Excel to keep storing them as text,
# matching the workbook's existing text-typed cells instead of converting them
# to real numbers/percentages.

$ws.Range("D2").Value = "'301.59"
$ws.Range("E2").Value = "'-0.20%"
$ws.Range("D3").Value = "'32.27"
$ws.Range("E3").Value = "'1.12%"
$ws.Range("D4").Value = "'5.000"
$ws.Range("E4").Value = "'-2.47%"
$ws.Range("D5").Value = "'0.07671"
$ws.Range("E5").Value = "'-2.30%"
$ws.Range("D6").Value = "'2.046"
$ws.Range("E6").Value = "'-10.46%"
$ws.Range("D7").Value = "'7.844"
$ws.Range("E7").Value = "'0.36%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'3.782"
$ws.Range("E8").Value = "'-0.75%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9158"
$ws.Range("E9").Value = "'-1.26%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1759"
$ws.Range("E10").Value = "'-1.11%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.07899"
$ws.Range("E11").Value = "'4.91%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08500"
$ws.Range("E12").Value = "'-5.40%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03073"
$ws.Range("E13").Value = "'-0.34%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09980"
$ws.Range("E14").Value = "'-0.40%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001511"
$ws.Range("E15").Value = "'-0.67%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005727"
$ws.Range("E16").Value = "'-2.38%"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "'0.007498"
$ws.Range("E17").Value = "'2,116.77%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.467"
$ws.Range("E18").Value = "'0.21%"
$ws.Range("D19").Value = "'2.153"
$ws.Range("E19").Value = "'-4.40%"
$ws.Range("E20").Value = "'1.42%"
$ws.Range("E21").Value = "'-0.83%"
$ws.Range("D22").Value = "'4.268"
$ws.Range("E22").Value = "'-1.79%"
$ws.Range("E23").Value = "'11.19%"
$ws.Range("D24").Value = "'0.04530"
$ws.Range("E24").Value = "'-1.74%"
$ws.Range("D25").Value = "'0.001229"
$ws.Range("E25").Value = "'-1.66%"
$ws.Range("D26").Value = "'0.004404"
$ws.Range("E26").Value = "'-1.78%"
$ws.Range("E27").Value = "'0.06%"
$ws.Range("D39").Value = "'0.01709"
$ws.Range("E39").Value = "'-4.14%"
$ws.Range("D40").Value = "'0.04673"
$ws.Range("E40").Value = "'-2.42%"
$ws.Range("D41").Value = "'0.007552"
$ws.Range("E41").Value = "'2.35%"
$ws.Range("D42").Value = "'0.1351"
$ws.Range("E42").Value = "'-0.94%"
$ws.Range("D43").Value = "'0.002330"
$ws.Range("E43").Value = "'6.46%"
$ws.Range("D44").Value = "'0.01052"
$ws.Range("E44").Value = "'8.24%"
$ws.Range("D45").Value = "'0.00006253"
$ws.Range("E45").Value = "'-0.74%"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("D47").Value = "'0.002999"
$ws.Range("E47").Value = "'-62.45%"
$ws.Range("E48").Value = "'14.59%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'0.05%"
